# Apply the edit: insert two new data rows (new rows 5 and 6) into the
# "Membrillo" sheet, pushing the existing rows 5..55 down to 7..57, and
# populate the two new rows with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows starting at row 5 (existing rows 5-55 shift to 7-57)
$ws.Rows.Item(5).Resize(2).Insert()

# ---- New row 5 ----
$ws.Cells.Item(5, 1).Value = 9
$ws.Cells.Item(5, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(5, 3).Value = "Metropolitana"
$ws.Cells.Item(5, 4).Value2 = 45022
$ws.Cells.Item(5, 5).Value = 13
$ws.Cells.Item(5, 6).Value = "Fruta"
$ws.Cells.Item(5, 7).Value = 100104
$ws.Cells.Item(5, 8).Value = "Frutos de pepita"
$ws.Cells.Item(5, 9).Value = 100104003
$ws.Cells.Item(5, 10).Value = "Membrillo"
$ws.Cells.Item(5, 11).Value = "Champion"
$ws.Cells.Item(5, 12).Value = "Especial"
$ws.Cells.Item(5, 13).Value = 270
$ws.Cells.Item(5, 14).Value = 12000
$ws.Cells.Item(5, 15).Value = 12000
$ws.Cells.Item(5, 16).Value = 12000
$ws.Cells.Item(5, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(5, 18).Value = "Región Metropolitana"
$ws.Cells.Item(5, 19).Value = 800
$ws.Cells.Item(5, 20).Value = 15

# ---- New row 6 ----
$ws.Cells.Item(6, 1).Value = 9
$ws.Cells.Item(6, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(6, 3).Value = "Metropolitana"
$ws.Cells.Item(6, 4).Value2 = 45022
$ws.Cells.Item(6, 5).Value = 13
$ws.Cells.Item(6, 6).Value = "Fruta"
$ws.Cells.Item(6, 7).Value = 100104
$ws.Cells.Item(6, 8).Value = "Frutos de pepita"
$ws.Cells.Item(6, 9).Value = 100104003
$ws.Cells.Item(6, 10).Value = "Membrillo"
$ws.Cells.Item(6, 11).Value = "Champion"
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 300
$ws.Cells.Item(6, 14).Value = 9000
$ws.Cells.Item(6, 15).Value = 9000
$ws.Cells.Item(6, 16).Value = 9000
$ws.Cells.Item(6, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(6, 18).Value = "Región Metropolitana"
$ws.Cells.Item(6, 19).Value = 600
$ws.Cells.Item(6, 20).Value = 15
